# SparkExecutorMemory.pptx - minor update
# Target shape: "TextBox 31" (the "User Memory: 1 -spark.memory.fraction -> 0.4" box)
# Changes:
#   1. Shrink the box width (and a 1-EMU nudge of its left edge that is below the
#      float precision of the Left/Width COM properties, see notes below).
#   2. Re-flow the second line's runs from
#        "1 -spark.memory.fraction -" + "> " + "0.4"
#      to
#        "1 " + "-spark.memory.fraction" + "=" + "0.4"

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(16)

# --- 1. Resize / reposition -------------------------------------------------
# Original EMU: off x=4664528 y=4546585, ext cx=1598677 cy=1477969
# Target   EMU: off x=4664529 y=4546585, ext cx=1540056 cy=1477969
# Shape.Left/.Top/.Width/.Height round-trip through a single-precision Single,
# so we feed the EMU-exact point values through; at this shape's magnitude the
# 1-EMU left nudge (4664528 -> 4664529) is below that precision floor (both
# collapse to the same Single), same as it would be in real PowerPoint COM
# automation, so Left is left untouched. Width moves by a clearly-resolvable
# amount, so it is set precisely (compensating by 1 EMU for the float
# round-trip so the saved value lands exactly on 1540056).
$sh.Left   = 4664529 / 12700.0
$sh.Top    = 4546585 / 12700.0
$sh.Width  = 1540057 / 12700.0
$sh.Height = 1477969 / 12700.0

# --- 2. Re-flow the second paragraph's runs ---------------------------------
$tr = $sh.TextFrame.TextRange

# Original 2nd paragraph = "1 -spark.memory.fraction -" + "> " + "0.4"
# (runs start at absolute offsets 14 / 40 / 42 in the full text range).

# a) Shrink the first run down to "1 " (keeps its original run properties).
$run1 = $tr.Characters(14, 26)
$run1.Text = "1 "

# b) The "> " run (now shifted left to offset 16) gets "-spark.memory.fraction="
#    typed in front of it - this lands in the same run as "> ", inheriting its
#    (smtClean) run properties.
$run2 = $tr.Characters(16, 2)
$run2.InsertBefore("-spark.memory.fraction=")

# c) Remove the now-unwanted "> " leftover (offset 39, length 2).
$leftover = $tr.Characters(39, 2)
$leftover.Delete()

# d) Remove the old "0.4" run (offset 39, length 3 after the above edits) and
#    retype it right after the "...fraction=" run so it picks up that run's
#    properties too.
$oldNum = $tr.Characters(39, 3)
$oldNum.Delete()
$fracEquals = $tr.Characters(16, 23)
$fracEquals.InsertAfter("0.4")
